$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.080.53"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.86%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.748.92"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.40%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9996"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.08%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").Style = "Normal"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5288"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.17%  "

$ws.Range("E8").Value = "  +0.59%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06187"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.68%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.743.51"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.30%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07184"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.89%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.37"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.38%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6447"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.87%  "

$ws.Range("E14").Value = "  +2.73%  "

$ws.Range("E15").Value = "  +2.52%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9998"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.02%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9995"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.15%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.009.08"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.51%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.68"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.16%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000006715"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.14%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.969.44"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.57%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.309"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.62%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.745"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.02%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.221"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.44%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "138.28"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.77%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.507"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.08%  "

$ws.Range("E27").Value = "  +2.31%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.802"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.62%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "104.57"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.86%  "

$ws.Range("E30").Value = "  +0.17%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.798"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.10%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.666"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.97%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04572"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.90%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.647"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.12%  "

$ws.Range("E35").Value = "  +3.77%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6332"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.02%  "

$ws.Range("E37").Value = "  +1.60%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01595"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.22%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.948"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.36%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.000"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.05%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "100.12"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.60%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.3926"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.95%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7452"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.43%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.028"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.09%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1144"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.03%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.342"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.67%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05350"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.18%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "30.92"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.37%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "54.10"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.84%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.675"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.63%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3454"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.16%  "
